# Update 5-Feb-2021, end of day update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 24: add 260000 to existing wages expense formula ---
$ws.Range("D24").Formula = "=60000+260000"

# --- Row 25: add 900000+600000 to the TRANSFER BCA formula ---
$ws.Range("D25").Formula = "=1550000+41600000+900000+600000"

# --- Row 26: add 26415000 to the A/R formula ---
$ws.Range("C26").Formula = "=17240000+24360000+26415000"

# --- Row 27: new entry, A/P ---
$ws.Range("B27").Value = "A/P"
$ws.Range("D27").Formula = "=1877000"

# --- Row 28: new entry, SALES - cash/retail ---
$ws.Range("B28").Value = "SALES - cash/retail"
$ws.Range("C28").Formula = "=4120475+30658525-26415000"

# --- Row 29: new entry, FREIGHT OUT ---
$ws.Range("B29").Value = "FREIGHT OUT"
$ws.Range("D29").Value = 54000

# --- Row 30: new entry, SELISIH - kurang ---
$ws.Range("B30").Value = "SELISIH - kurang"
$ws.Range("D30").Value = 425000

# --- Row 31: new entry, SETOR KE BANK ---
$ws.Range("B31").Value = "SETOR KE BANK"
$ws.Range("D31").Formula = "=29000000"

# --- Row 32: new entry, DOKTER - qiu ---
$ws.Range("B32").Value = "DOKTER - qiu"
$ws.Range("D32").Value = 800000

# --- Row 33: new day (5-Feb-2021), Wages Expense ---
$ws.Range("A33").Value = 44232
$ws.Range("B33").Value = "Wages Expense"
$ws.Range("D33").Formula = "=60000"

# --- Row 34: TRANSFER BCA ---
$ws.Range("B34").Value = "TRANSFER BCA"
$ws.Range("D34").Formula = "=2216000+154800000+1130000+55000"

# --- Row 35: A/R ---
$ws.Range("B35").Value = "A/R"
$ws.Range("C35").Formula = "=154800000"

# --- Update the frozen pane / active selection to reflect the new view state ---
$win = $excel.ActiveWindow
$win.ScrollRow = 32
$ws.Range("D53").Select()
